# Adds a new "2022-Q1" sheet (holdings detail) positioned right before the
# "总计" (totals) summary sheet, and inserts a corresponding summary row at
# the top of the "总计" sheet's data.
#
# NOTE: worksheet variables in this runtime are bound by position, not by
# identity, so after any operation that reorders sheets (Add, Move, ...)
# we must re-fetch sheet references (by name) rather than reuse old ones.

function Set-TextCell($ws, $rowIdx, $colIdx, $text) {
    $cell = $ws.Cells.Item($rowIdx, $colIdx)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" worksheet right before "总计".
# ---------------------------------------------------------------------
$totalSheetOld = $wb.Worksheets.Item($wb.Worksheets.Count)
$q1Sheet = $wb.Worksheets.Add($totalSheetOld)
$q1Sheet.Name = "2022-Q1"

# Reuse formatting (bold/centered/bordered style) from an existing,
# identically-shaped sheet so no redundant styles are introduced.
$template = $wb.Worksheets.Item("2021-Q4")
$template.Range("B1:H1").Copy()
$q1Sheet.Range("B1").PasteSpecial(-4122)   # xlPasteFormats
$template.Range("A2:A12").Copy()
$q1Sheet.Range("A2").PasteSpecial(-4122)   # xlPasteFormats

# Header row
$q1Sheet.Cells.Item(1,2).Value = "基金代码"
$q1Sheet.Cells.Item(1,3).Value = "基金名称"
$q1Sheet.Cells.Item(1,4).Value = "基金规模"
$q1Sheet.Cells.Item(1,5).Value = "股票总仓位"
$q1Sheet.Cells.Item(1,6).Value = "仓位占比"
$q1Sheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$q1Sheet.Cells.Item(1,8).Value = "仓位排名"

# Data rows: index, 基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名
$rows = @(
    @(0, "003230", "创金合信医疗保健行业股票A",                 "8.07",  "94.55", "7.38", "0.5956", 6),
    @(1, "009960", "银华多元机遇混合",                           "10.75", "87.53", "4.24", "0.4558", 3),
    @(2, "003231", "创金合信医疗保健行业股票C",                 "4.28",  "94.55", "7.38", "0.3159", 6),
    @(3, "010585", "创金合信医药消费股票A",                     "4.22",  "93.65", "5.22", "0.2203", 6),
    @(4, "009246", "摩根士丹利华鑫ESG量化先行混合",             "4.04",  "92.17", "2.22", "0.0897", 3),
    @(5, "013067", "富安达中小盘六个月持有期混合",               "2.45",  "74.39", "3.36", "0.0823", 5),
    @(6, "011383", "富安达医药创新混合",                         "1.68",  "83.50", "4.13", "0.0694", 4),
    @(7, "180028", "银华永祥灵活配置混合",                       "0.61",  "77.23", "5.25", "0.0320", 3),
    @(8, "001861", "富安达健康人生灵活配置混合",                 "0.61",  "82.18", "4.63", "0.0282", 5),
    @(9, "010586", "创金合信医药消费股票C",                     "0.50",  "93.65", "5.22", "0.0261", 6),
    @(10, "004536", "嘉实中小企业量化活力灵活配置混合",         "0.17",  "90.06", "1.48", "0.0025", 8)
)

$r = 2
foreach ($row in $rows) {
    $q1Sheet.Cells.Item($r, 1).Value = $row[0]
    Set-TextCell $q1Sheet $r 2 $row[1]
    $q1Sheet.Cells.Item($r, 3).Value = $row[2]
    Set-TextCell $q1Sheet $r 4 $row[3]
    Set-TextCell $q1Sheet $r 5 $row[4]
    Set-TextCell $q1Sheet $r 6 $row[5]
    Set-TextCell $q1Sheet $r 7 $row[6]
    $q1Sheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Insert the "2022-Q1" summary row at the top of "总计" data, pushing
#    the existing rows down by one. Re-fetch "总计" by name since the
#    sheet collection order changed in step 1.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Cells.Item(2,1).Value = 0
$totalSheet.Cells.Item(2,2).Value = "2022-Q1"
$totalSheet.Cells.Item(2,3).Value = 11
$totalSheet.Cells.Item(2,4).Value = 1.92

# Re-apply the column-A style (Insert sometimes leaves the new row's
# other cells without explicit formatting) and renumber the index
# column sequentially for the rows that shifted down.
$totalSheet.Cells.Item(2,1).NumberFormat = $totalSheet.Cells.Item(3,1).NumberFormat

for ($i = 3; $i -le 7; $i++) {
    $totalSheet.Cells.Item($i, 1).Value = $i - 2
}

Write-Host "done"
